$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet updates (publish a refreshed FHIR logical model)
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value  = "2.0.0-sd-202406-matchbox-patch"
$meta.Range("B8").Value  = "2024-06-19T17:47:42+02:00"
$meta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# ---------------------------------------------------------------------------
# 2. Elements sheet: insert a new row for CustodianOrganization.sdtcTelecom
#    right above the existing CustodianOrganization.addr row (row 17),
#    pushing addr down to row 18.
# ---------------------------------------------------------------------------
$els = $wb.Worksheets.Item("Elements")

# Shift row 17 (addr) and below down by one row, making room for the new row.
$els.Rows.Item(17).Insert()

# The new "sdtcTelecom" element row has exactly the same shape/cardinality as
# the existing "telecom" row just above it (row 16: Min 0, Max *, type TEL) -
# only the element path differs. Clone row 16's formatting and values onto
# the freshly inserted row 17, then fix up the path-name cells.
$els.Range("A16:AK16").Copy()
$els.Range("A17:AK17").PasteSpecial(-4122)
$els.Range("A16:AK16").Copy()
$els.Range("A17:AK17").PasteSpecial(-4163)

$els.Range("A17").Value  = "CustodianOrganization.sdtcTelecom"
$els.Range("B17").Value  = "CustodianOrganization.sdtcTelecom"
$els.Range("AF17").Value = "CustodianOrganization.sdtcTelecom"
